$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 239.2

$ws.Range("B3").Value = 228.2
$ws.Range("C3").Value = 205.4

$ws.Range("C4").Value = 103.8

$ws.Range("C5").Value = -27.3

$ws.Range("C6").Value = -19.7

$ws.Range("C7").Value = -147.1

$ws.Range("C8").Value = -162.6

$ws.Range("C9").Value = -130

$ws.Range("C10").Value = -82.3

$ws.Range("C11").Value = 45.6

$ws.Range("C12").Value = 129.4

$ws.Range("C13").Value = 65.5

$ws.Range("C14").Value = 219.8

$ws.Range("C15").Value = 196

$ws.Range("C16").Value = 91.09999999999999

$ws.Range("C17").Value = 57.2

$ws.Range("C18").Value = 171.1

$ws.Range("C19").Value = 101.9

$ws.Range("C21").Value = 264.8

$ws.Range("C22").Value = 398.4

$ws.Range("C23").Value = 522.7

$ws.Range("C24").Value = 275.3
